$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A8").Value = "2025-1"
$ws.Range("B8").Value = "EP TASA 418"
$ws.Range("C8").Value = "Embarcación Pesquera"
$ws.Range("D8").Value = "GP/80"
$ws.Range("E8").Value = "GP/80-125"

$ws.Range("A9").Value = "2025-1"
$ws.Range("B9").Value = "EP TASA 412"
$ws.Range("C9").Value = "Embarcación Pesquera"
$ws.Range("D9").Value = "GP/79"
$ws.Range("E9").Value = "GP/79-125"

$ws.Range("F13").Select()
